$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '23.209.86'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '1.603.58'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '0.9998'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").Value = '303.35'
$ws.Range("D7").Value = '0.3783'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '51.88'
$ws.Range("E8").Value = '  +2.88%  '
$ws.Range("D9").Value = '0.3628'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '1.271'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").Value = '0.9998'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '0.08122'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '7.422'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '0.00001247'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '1.602.81'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '93.95'
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("D19").Value = '0.06873'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '12.96'
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").Value = '23.203.61'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '3.012'
$ws.Range("E25").Value = '  +8.05%  '
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").Value = '21.23'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").Value = '150.00'
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '5.247'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").Value = '133.91'
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").Value = '2.370'
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '6.754'
$ws.Range("E32").Value = '  -1.28%  '
$ws.Range("D33").Value = '1.779.01'
$ws.Range("D34").Value = '0.9665'
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").Value = '0.07518'
$ws.Range("E35").Value = '  -2.24%  '
$ws.Range("D36").Value = '0.02727'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '10.24'
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").Value = '0.2527'
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").Value = '6.086'
$ws.Range("E40").Value = '  -3.21%  '
$ws.Range("D41").Value = '1.369'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '0.7107'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = '12.55'
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("D44").Value = '15.59'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D46").Value = '2.316'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").Value = '4.021'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").Value = '132.34'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = '0.07955'
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("E51").Value = '  +0.80%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
